# Generate Report for Handoff
# Adds a new row (row 9) describing the file
# "b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md" to all three report sheets
# (Overview, zh-cn, de-de) and grows the backing tables / ranges to match.

$wb = $excel.ActiveWorkbook

# A single apostrophe forces the cell to be treated as literal text; once
# the leading quote marker is stripped what remains is an empty string, so
# this is how we write an explicit (non-blank) empty-string cell instead of
# leaving the cell absent altogether - matching the existing sheets, which
# use an empty shared string for "no value" fields.
$EMPTY = "'"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f7f7c2a6f9f7f0e6c4f0a0a2b4c6d8e0f1a2b3c/e2e/b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md"

# ---------------------------------------------------------------------
# Sheet "Overview" - row 9
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(9, 1).Value = "b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md"
$ov.Cells.Item(9, 2).Value = "e2e\b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md"
$ov.Hyperlinks.Add($ov.Cells.Item(9, 2), $baseUrl, "", "", "e2e\b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md")
$ov.Cells.Item(9, 3).Value = ".md"
$ov.Cells.Item(9, 4).Value = $EMPTY
$ov.Cells.Item(9, 5).Value = "Ready for handoff"
$ov.Cells.Item(9, 6).Value = "Ready for handoff"
$ov.Cells.Item(9, 7).Value = "2016-08-25 22:43:20"
$ov.Cells.Item(9, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G9"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" - row 9
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Cells.Item(9, 1).Value = "b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md"
$zh.Hyperlinks.Add($zh.Cells.Item(9, 1), $baseUrl, "", "", "b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md")
$zh.Cells.Item(9, 2).Value = ".md"
$zh.Cells.Item(9, 3).Value = "Ready for handoff"
$zh.Cells.Item(9, 4).Value = "e2e"
$zh.Cells.Item(9, 5).Value = "ht"
$zh.Cells.Item(9, 6).Value = "False"
$zh.Cells.Item(9, 7).Value = "b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.46ae93fd9c3ce5ceaf0c79ab66125ce0a1aa81e1.zh-cn.xlf"
$zh.Cells.Item(9, 8).Value = "2016-08-25 22:43:15"
$zh.Cells.Item(9, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Cells.Item(9, 9).Value = $EMPTY
$zh.Cells.Item(9, 10).Value = $EMPTY
$zh.Cells.Item(9, 11).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(9, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Cells.Item(9, 12).Value = $EMPTY
$zh.Cells.Item(9, 13).Value = "True"
$zh.Cells.Item(9, 14).Value = $EMPTY
$zh.Cells.Item(9, 15).Value = "False"
$zh.Cells.Item(9, 16).Value = $EMPTY

$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:P9"))

# ---------------------------------------------------------------------
# Sheet "de-de" - row 9
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Cells.Item(9, 1).Value = "b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md"
$de.Hyperlinks.Add($de.Cells.Item(9, 1), $baseUrl, "", "", "b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md")
$de.Cells.Item(9, 2).Value = ".md"
$de.Cells.Item(9, 3).Value = "Ready for handoff"
$de.Cells.Item(9, 4).Value = "e2e"
$de.Cells.Item(9, 5).Value = "ht"
$de.Cells.Item(9, 6).Value = "False"
$de.Cells.Item(9, 7).Value = "b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.46ae93fd9c3ce5ceaf0c79ab66125ce0a1aa81e1.de-de.xlf"
$de.Cells.Item(9, 8).Value = "2016-08-25 22:43:20"
$de.Cells.Item(9, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Cells.Item(9, 9).Value = $EMPTY
$de.Cells.Item(9, 10).Value = $EMPTY
$de.Cells.Item(9, 11).Value = "0001-01-01 00:00:00"
$de.Cells.Item(9, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Cells.Item(9, 12).Value = $EMPTY
$de.Cells.Item(9, 13).Value = "True"
$de.Cells.Item(9, 14).Value = $EMPTY
$de.Cells.Item(9, 15).Value = "False"
$de.Cells.Item(9, 16).Value = $EMPTY

$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:P9"))

Write-Host "Added handoff row for b135cdc8-9f0c-4e91-a0bf-c8f5eae628e4.md to Overview, zh-cn, de-de."
